$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5328953266143799
$ws.Range("B1").Value = 0.6724743247032166
$ws.Range("C1").Value = 5.751636505126953
$ws.Range("D1").Value = 1.56080424785614
$ws.Range("E1").Value = 1.002545952796936
